$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema3a"
$ws.Range("C2").Value = "Plxna2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8200746666666667
$ws.Range("H2").Value = 2.460224
$ws.Range("I2").Value = 0.04311293902675128
$ws.Range("J2").Value = 0.04311293902675128
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 40.25420133333333
$ws.Range("N2").Value = 120.762604
$ws.Range("O2").Value = 0.4854671023051697
$ws.Range("P2").Value = 0.4854671023051695
$ws.Range("Q2").Value = 33.01145074036622
$ws.Range("R2").Value = 297.103056663296
$ws.Range("S2").Value = 0.02092991358117641
$ws.Range("T2").Value = 0.0209299135811764

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema3a"
$ws.Range("C3").Value = "Plxna2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8200746666666667
$ws.Range("H3").Value = 2.460224
$ws.Range("I3").Value = 0.04311293902675128
$ws.Range("J3").Value = 0.04311293902675128
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.624984
$ws.Range("N3").Value = 13.874952
$ws.Range("O3").Value = 0.05577747182450057
$ws.Range("P3").Value = 0.05577747182450056
$ws.Range("Q3").Value = 3.792832212138667
$ws.Range("R3").Value = 34.135489909248
$ws.Range("S3").Value = 0.002404730741836031
$ws.Range("T3").Value = 0.00240473074183603

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema3a"
$ws.Range("C4").Value = "Plxna2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8200746666666667
$ws.Range("H4").Value = 2.460224
$ws.Range("I4").Value = 0.04311293902675128
$ws.Range("J4").Value = 0.04311293902675128
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 36.68940733333334
$ws.Range("N4").Value = 110.068222
$ws.Range("O4").Value = 0.4424755596543956
$ws.Range("P4").Value = 0.4424755596543954
$ws.Range("Q4").Value = 30.0880534890809
$ws.Range("R4").Value = 270.7924814017281
$ws.Range("S4").Value = 0.0190764218242076
$ws.Range("T4").Value = 0.0190764218242076

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sema3a"
$ws.Range("C5").Value = "Plxna2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.8200746666666667
$ws.Range("H5").Value = 2.460224
$ws.Range("I5").Value = 0.04311293902675128
$ws.Range("J5").Value = 0.04311293902675128
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.349902
$ws.Range("N5").Value = 4.049706
$ws.Range("O5").Value = 0.01627986621593436
$ws.Range("P5").Value = 0.01627986621593436
$ws.Range("Q5").Value = 1.107020432682667
$ws.Range("R5").Value = 9.963183894144002
$ws.Range("S5").Value = 0.0007018728795312463
$ws.Range("T5").Value = 0.0007018728795312461

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema3a"
$ws.Range("C6").Value = "Plxna2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.358031333333334
$ws.Range("H6").Value = 4.074094000000001
$ws.Range("I6").Value = 0.07139437962203982
$ws.Range("J6").Value = 0.07139437962203982
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 40.25420133333333
$ws.Range("N6").Value = 120.762604
$ws.Range("O6").Value = 0.4854671023051697
$ws.Range("P6").Value = 0.4854671023051695
$ws.Range("Q6").Value = 54.66646670897511
$ws.Range("R6").Value = 491.998200380776
$ws.Range("S6").Value = 0.03465962259598692
$ws.Range("T6").Value = 0.03465962259598691

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema3a"
$ws.Range("C7").Value = "Plxna2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.358031333333334
$ws.Range("H7").Value = 4.074094000000001
$ws.Range("I7").Value = 0.07139437962203982
$ws.Range("J7").Value = 0.07139437962203982
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.624984
$ws.Range("N7").Value = 13.874952
$ws.Range("O7").Value = 0.05577747182450057
$ws.Range("P7").Value = 0.05577747182450056
$ws.Range("Q7").Value = 6.280873188165335
$ws.Range("R7").Value = 56.52785869348801
$ws.Range("S7").Value = 0.003982197997796024
$ws.Range("T7").Value = 0.003982197997796023

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sema3a"
$ws.Range("C8").Value = "Plxna2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.358031333333334
$ws.Range("H8").Value = 4.074094000000001
$ws.Range("I8").Value = 0.07139437962203982
$ws.Range("J8").Value = 0.07139437962203982
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 36.68940733333334
$ws.Range("N8").Value = 110.068222
$ws.Range("O8").Value = 0.4424755596543956
$ws.Range("P8").Value = 0.4424755596543954
$ws.Range("Q8").Value = 49.82536476009646
$ws.Range("R8").Value = 448.4282828408681
$ws.Range("S8").Value = 0.03159026807944044
$ws.Range("T8").Value = 0.03159026807944043

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sema3a"
$ws.Range("C9").Value = "Plxna2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.358031333333334
$ws.Range("H9").Value = 4.074094000000001
$ws.Range("I9").Value = 0.07139437962203982
$ws.Range("J9").Value = 0.07139437962203982
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.349902
$ws.Range("N9").Value = 4.049706
$ws.Range("O9").Value = 0.01627986621593436
$ws.Range("P9").Value = 0.01627986621593436
$ws.Range("Q9").Value = 1.833209212929334
$ws.Range("R9").Value = 16.498882916364
$ws.Range("S9").Value = 0.001162290948816439
$ws.Range("T9").Value = 0.001162290948816438

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Sema3a"
$ws.Range("C10").Value = "Plxna2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 16.8273
$ws.Range("H10").Value = 50.4819
$ws.Range("I10").Value = 0.8846442749337277
$ws.Range("J10").Value = 0.8846442749337278
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 40.25420133333333
$ws.Range("N10").Value = 120.762604
$ws.Range("O10").Value = 0.4854671023051697
$ws.Range("P10").Value = 0.4854671023051695
$ws.Range("Q10").Value = 677.3695220963998
$ws.Range("R10").Value = 6096.325698867598
$ws.Range("S10").Value = 0.4294656927229346
$ws.Range("T10").Value = 0.4294656927229346

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Sema3a"
$ws.Range("C11").Value = "Plxna2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 16.8273
$ws.Range("H11").Value = 50.4819
$ws.Range("I11").Value = 0.8846442749337277
$ws.Range("J11").Value = 0.8846442749337278
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.624984
$ws.Range("N11").Value = 13.874952
$ws.Range("O11").Value = 0.05577747182450057
$ws.Range("P11").Value = 0.05577747182450056
$ws.Range("Q11").Value = 77.82599326319999
$ws.Range("R11").Value = 700.4339393688
$ws.Range("S11").Value = 0.04934322111982173
$ws.Range("T11").Value = 0.04934322111982173

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Sema3a"
$ws.Range("C12").Value = "Plxna2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 16.8273
$ws.Range("H12").Value = 50.4819
$ws.Range("I12").Value = 0.8846442749337277
$ws.Range("J12").Value = 0.8846442749337278
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 36.68940733333334
$ws.Range("N12").Value = 110.068222
$ws.Range("O12").Value = 0.4424755596543956
$ws.Range("P12").Value = 0.4424755596543954
$ws.Range("Q12").Value = 617.3836640202001
$ws.Range("R12").Value = 5556.4529761818
$ws.Range("S12").Value = 0.3914334706463581
$ws.Range("T12").Value = 0.3914334706463581

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Sema3a"
$ws.Range("C13").Value = "Plxna2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 16.8273
$ws.Range("H13").Value = 50.4819
$ws.Range("I13").Value = 0.8846442749337277
$ws.Range("J13").Value = 0.8846442749337278
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.349902
$ws.Range("N13").Value = 4.049706
$ws.Range("O13").Value = 0.01627986621593436
$ws.Range("P13").Value = 0.01627986621593436
$ws.Range("Q13").Value = 22.7152059246
$ws.Range("R13").Value = 204.4368533214
$ws.Range("S13").Value = 0.01440189044461334
$ws.Range("T13").Value = 0.01440189044461334

$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Sema3a"
$ws.Range("C14").Value = "Plxna2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.016138
$ws.Range("H14").Value = 0.048414
$ws.Range("I14").Value = 0.0008484064174811467
$ws.Range("J14").Value = 0.0008484064174811467
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 40.25420133333333
$ws.Range("N14").Value = 120.762604
$ws.Range("O14").Value = 0.4854671023051697
$ws.Range("P14").Value = 0.4854671023051695
$ws.Range("Q14").Value = 0.6496223011173332
$ws.Range("R14").Value = 5.846600710055999
$ws.Range("S14").Value = 0.0004118734050716823
$ws.Range("T14").Value = 0.0004118734050716822

$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Sema3a"
$ws.Range("C15").Value = "Plxna2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.016138
$ws.Range("H15").Value = 0.048414
$ws.Range("I15").Value = 0.0008484064174811467
$ws.Range("J15").Value = 0.0008484064174811467
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.624984
$ws.Range("N15").Value = 13.874952
$ws.Range("O15").Value = 0.05577747182450057
$ws.Range("P15").Value = 0.05577747182450056
$ws.Range("Q15").Value = 0.07463799179200001
$ws.Range("R15").Value = 0.671741926128
$ws.Range("S15").Value = 0.000047321965046780121457413448
$ws.Range("T15").Value = 0.000047321965046780107904886292

$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Sema3a"
$ws.Range("C16").Value = "Plxna2"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.016138
$ws.Range("H16").Value = 0.048414
$ws.Range("I16").Value = 0.0008484064174811467
$ws.Range("J16").Value = 0.0008484064174811467
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 36.68940733333334
$ws.Range("N16").Value = 110.068222
$ws.Range("O16").Value = 0.4424755596543956
$ws.Range("P16").Value = 0.4424755596543954
$ws.Range("Q16").Value = 0.5920936555453334
$ws.Range("R16").Value = 5.328842899908
$ws.Range("S16").Value = 0.0003753991043893511
$ws.Range("T16").Value = 0.0003753991043893511

$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Sema3a"
$ws.Range("C17").Value = "Plxna2"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.016138
$ws.Range("H17").Value = 0.048414
$ws.Range("I17").Value = 0.0008484064174811467
$ws.Range("J17").Value = 0.0008484064174811467
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.349902
$ws.Range("N17").Value = 4.049706
$ws.Range("O17").Value = 0.01627986621593436
$ws.Range("P17").Value = 0.01627986621593436
$ws.Range("Q17").Value = 0.021784718476
$ws.Range("R17").Value = 0.196062466284
$ws.Range("S17").Value = 0.000013811942973333220837772045
$ws.Range("T17").Value = 0.000013811942973333220837772045
